$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

$ws.Cells.Item($row, 1).Value = "2025-05-01T11:50:26.484Z"
$ws.Cells.Item($row, 2).Value = "NRC"
$ws.Cells.Item($row, 3).Value = "C3"
$ws.Cells.Item($row, 4).Value = "الرحلة 2"
$ws.Cells.Item($row, 5).Value = "ايتا"
$ws.Cells.Item($row, 6).Value = "احمد"
$ws.Cells.Item($row, 7).Value = "2323"
$ws.Cells.Item($row, 8).Value = ""
